$d = $word.ActiveDocument

# --- Paragraphs 1-3: replace placeholder "---" lines with names/net IDs ---
$d.Paragraphs.Item(1).Range.Text = "Adarsh Raghupati   axh190002"
$d.Paragraphs.Item(2).Range.Text = "Akash Akki         apa190001"
$d.Paragraphs.Item(3).Range.Text = "Keerti Keerti      kxk190012"

# --- Paragraph 4: "Stewart Cannon (sjc160330)" -> split across three runs ---
$p4 = $d.Paragraphs.Item(4)
$p4xml = '<w:p w14:paraId="275BE21D" w14:textId="6B9BDC27" w:rsidR="008E0B30" w:rsidRDefault="008E0B30" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:r><w:t xml:space="preserve">Stewart </w:t></w:r><w:r><w:t>C</w:t></w:r><w:r><w:t>annon     sjc160330</w:t></w:r></w:p>'
$p4.Range.InsertXML($p4xml)

# --- Paragraph 7: skip list paragraph -> append memory-use sentence ---
$p7 = $d.Paragraphs.Item(7)
$ip7 = $d.Range($p7.Range.End - 1, $p7.Range.End - 1)
$ip7.InsertAfter(" Memory use for the skip list was also generally between 45% and 65% over the tests.")

# --- Paragraph 8: red black tree paragraph -> append memory-use sentence ---
$p8 = $d.Paragraphs.Item(8)
$ip8 = $d.Range($p8.Range.End - 1, $p8.Range.End - 1)
$ip8.InsertAfter(" The memory for the red black tree also fluctuated between 45% and 78% memory use, with two values being close at 73% and 78%.")

# --- Paragraph 9: TreeSet paragraph -> append memory-use sentence ---
$p9 = $d.Paragraphs.Item(9)
$ip9 = $d.Range($p9.Range.End - 1, $p9.Range.End - 1)
$ip9.InsertAfter(" The memory usage of TreeSet showed a higher utilization of 85% and 65% for 4 and 8 million, while the utilization for 16 and 32 million was only 31% and 39%.")

# --- Paragraph 10: conclusion paragraph -> append memory-use sentence BEFORE the _GoBack bookmark ---
$p10 = $d.Paragraphs.Item(10)
$p10xml = '<w:p w14:paraId="0BA6AB6A" w14:textId="33862356" w:rsidR="005D3FF1" w:rsidRDefault="005D3FF1" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:r><w:tab/><w:t xml:space="preserve">In conclusion, </w:t></w:r><w:r w:rsidR="00D76F6B"><w:t xml:space="preserve">the </w:t></w:r><w:r w:rsidR="00BE19C0"><w:t xml:space="preserve">skip list was clearly the slowest, while the red black tree and the TreeSet </w:t></w:r><w:r w:rsidR="00A65174"><w:t>were more competitive</w:t></w:r><w:r w:rsidR="008677EC"><w:t xml:space="preserve">. However, the TreeSet structure was </w:t></w:r><w:r w:rsidR="00A820AC"><w:t>still the fastest of the three.</w:t></w:r><w:r><w:t xml:space="preserve"> Looking at memory usage, the values seen seemed unstable, however the TreeSet data structure consistently used less amounts of memory than the other structures at higher numbers of elements.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$p10.Range.InsertXML($p10xml)
